$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank column at N (14th column), shifting existing N:S to O:T
$ws.Columns("N:N").Insert()

# Match the column width of the neighboring "object" columns (9-13)
$ws.Columns("N:N").ColumnWidth = 17.5

# Set header for the new "subject" column
$ws.Range("N1").Value = "subject"

# Fill in values for the rows that get a subject value
$ws.Range("N4").Value = ":Person"
$ws.Range("N8").Value = ":Person"
$ws.Range("N20").Value = ":Image"

# Set the selection to A20 as shown in the diff
$ws.Range("A20").Select()
